$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2023-08-28 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-29 Tuesday", 2)

# Update the division problems in the table. Each cell is addressed
# individually (Tables.Item(1).Cell(row, col)) so that replacements never
# bleed into each other, even when a "new" value happens to equal another
# cell's "old" value (e.g. 83÷5= -> 26÷5=, while a different cell's
# 26÷5= -> 84÷5=).
$t = $d.Tables.Item(1)

function Replace-CellText($table, $row, $col, $old, $new) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Row 1
Replace-CellText $t 1 1 "47÷4=" "76÷7="
Replace-CellText $t 1 2 "56÷2=" "22÷7="
Replace-CellText $t 1 3 "43÷9=" "39÷7="
Replace-CellText $t 1 4 "92÷8=" "12÷4="
Replace-CellText $t 1 5 "89÷4=" "97÷8="

# Row 5
Replace-CellText $t 5 1 "58÷6=" "42÷2="
Replace-CellText $t 5 2 "34÷5=" "83÷8="
Replace-CellText $t 5 3 "65÷2=" "67÷3="
Replace-CellText $t 5 4 "34÷7=" "63÷3="
Replace-CellText $t 5 5 "64÷5=" "25÷9="

# Row 9
Replace-CellText $t 9 1 "26÷5=" "84÷5="
Replace-CellText $t 9 2 "37÷8=" "63÷6="
Replace-CellText $t 9 3 "13÷7=" "52÷7="
Replace-CellText $t 9 4 "86÷3=" "44÷2="
Replace-CellText $t 9 5 "93÷7=" "35÷2="

# Row 13
Replace-CellText $t 13 1 "77÷3=" "49÷6="
Replace-CellText $t 13 2 "71÷2=" "55÷8="
Replace-CellText $t 13 3 "32÷6=" "23÷4="
Replace-CellText $t 13 4 "83÷5=" "26÷5="
Replace-CellText $t 13 5 "46÷8=" "84÷3="

# Row 17
Replace-CellText $t 17 1 "10÷9=" "99÷5="
Replace-CellText $t 17 2 "73÷2=" "19÷8="
Replace-CellText $t 17 3 "48÷8=" "60÷7="
Replace-CellText $t 17 4 "84÷8=" "80÷7="
Replace-CellText $t 17 5 "45÷9=" "31÷3="
